$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Relocate the "_GoBack" bookmark from right after "$ touch FILE" to
# a point inside "$ Ssh-keygen " (right after the first letter), and
# fix the capitalisation typo "Ssh-keygen" -> "ssh-keygen" at the
# same time - without disturbing the separate "-" / "t rsa" runs
# that follow in that paragraph.
# ------------------------------------------------------------------

# Locate the word "Ssh-keygen".
$findRng = $d.Content
$findRng.Find.Execute("Ssh-keygen", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
if (-not $findRng.Find.Found) {
    throw "Could not locate 'Ssh-keygen' in the document."
}

$sPos = $findRng.Start        # position of the leading capital "S"
$splitPos = $sPos + 1         # position right after that letter

# Temporary bookmark right before the "S" - keeps the upcoming text
# edit from merging that run together with the preceding "$ " run.
$d.Bookmarks.Add("ZZTMP_SPLIT", $d.Range($sPos, $sPos))

# Re-adding "_GoBack" at the new location silently relocates it (Word
# bookmark names are unique, so this also removes the old occurrence
# that used to sit right after "$ touch FILE"). Placing it here first
# also keeps the text edit below from merging into "sh-keygen ".
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))

# Fix the typo: capital "S" -> lowercase "s".
$letterRng = $d.Range($sPos, $splitPos)
$letterRng.Text = "s"

# Drop the temporary helper bookmark - it was only scaffolding.
$d.Bookmarks("ZZTMP_SPLIT").Delete()
